$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NAText {
    param($cellRef, $textVal, $refStyleCell)
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $textVal
    $ws.Range($refStyleCell).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

# --- Simple numeric value updates ---
$ws.Range("G15").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("I16").Value = 9
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = -25
$ws.Range("L16").Value = 0
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -26.086956521739
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = -35.294117647058
$ws.Range("L17").Value = -8.333333333333
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 50
$ws.Range("L18").Value = -36.363636363636
$ws.Range("C19").Value = 9
$ws.Range("E19").Value = 80
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 64.705882352941
$ws.Range("I19").Value = 48
$ws.Range("J19").Value = 39
$ws.Range("K19").Value = 23.076923076923
$ws.Range("L19").Value = -27.272727272727
$ws.Range("D20").Value = 4
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -80
$ws.Range("J20").Value = 11
$ws.Range("K20").Value = -72.727272727272
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 55
$ws.Range("H21").Value = -3.508771929824
$ws.Range("I21").Value = 89
$ws.Range("J21").Value = 104
$ws.Range("K21").Value = -14.423076923076
$ws.Range("L21").Value = -29.365079365079
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -100
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = -66.666666666666
$ws.Range("L23").Value = -66.666666666666
$ws.Range("C24").Value = 29
$ws.Range("E24").Value = 20.833333333333
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 41.573033707865
$ws.Range("I24").Value = 215
$ws.Range("J24").Value = 166
$ws.Range("K24").Value = 29.518072289156
$ws.Range("L24").Value = 35.220125786163
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 2.5
$ws.Range("I25").Value = 68
$ws.Range("J25").Value = 68
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 23.636363636363
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -100
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 60
$ws.Range("I30").Value = 3

# --- N/A (text placeholder) -> numeric ---
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E23").Value = -100
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E28").Value = -100
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("G28").Value = 1
$ws.Range("H28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H28").Value = -100
$ws.Range("J28").NumberFormat = "#,##0"
$ws.Range("J28").Value = 1
$ws.Range("K28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K28").Value = -100
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E29").Value = -100
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("G29").Value = 1
$ws.Range("H29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H29").Value = -100
$ws.Range("J29").NumberFormat = "#,##0"
$ws.Range("J29").Value = 1
$ws.Range("K29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K29").Value = -100
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("C30").Value = 1

# --- numeric -> N/A (text placeholder) ---
Set-NAText "F23" "0" "F22"
Set-NAText "F26" "0" "F22"
Set-NAText "D27" "0" "D22"
Set-NAText "E27" "***.*" "E22"

$excel.CutCopyMode = $false

# --- Shared-string rich text edits (header) ---
$ws.Range("A8").Characters(21, 1).Text = "7"
$ws.Range("C9").Characters(46, 9).Text = "2/18/2024"
$ws.Range("C9").Characters(27, 8).Text = "2/12/2024"
